$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 currently holds the values as text ("4141.0" / "4142.0").
# Convert them to real numbers 4141 / 4142.
$ws.Range("A6").Value = 4141
$ws.Range("B6").Value = 4142

# Insert a new row 7 with text values "4000.0" / "4000.0".
# Force text storage (so Excel doesn't auto-convert the numeric-looking
# string into a number), then clear the temporary formatting so the
# cells end up with no special style, matching a plain text entry.
$ws.Range("A7:B7").NumberFormat = "@"
$ws.Range("A7").Value = "4000.0"
$ws.Range("B7").Value = "4000.0"
$ws.Range("A7:B7").ClearFormats()
